$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 963
$ws.Range("I6").Value = 963
$ws.Range("K6").Value = 2889
$ws.Range("M6").Value = -2777

$ws.Range("H76").Value = 20002646
$ws.Range("I76").Value = 28574068
$ws.Range("J76").Value = 2659.4666
$ws.Range("K76").Value = 28574068
$ws.Range("L76").Value = 2659.4666
$ws.Range("M76").Value = -28573753
$ws.Range("N76").Value = -3289.4666

$ws.Range("H79").Value = 20002646
$ws.Range("I79").Value = 28574068
$ws.Range("J79").Value = 2659.4666
$ws.Range("K79").Value = 28574068
$ws.Range("L79").Value = 2659.4666
$ws.Range("M79").Value = -28572976
$ws.Range("N79").Value = -4843.4666

$ws.Range("H81").Value = 27980
$ws.Range("J81").Value = 27980
$ws.Range("L81").Value = 27980
$ws.Range("N81").Value = -29976

$ws.Range("H84").Value = 27980
$ws.Range("J84").Value = 27980
$ws.Range("L84").Value = 83940
$ws.Range("N84").Value = -93924

$ws.Range("H116").Value = 9338654
$ws.Range("I116").Value = 4168960.8
$ws.Range("J116").Value = 18529220
$ws.Range("K116").Value = 4168960.8
$ws.Range("L116").Value = 18529220
$ws.Range("M116").Value = -4165518.8
$ws.Range("N116").Value = -18536104

$ws.Range("H138").Value = 2487.3186
$ws.Range("I138").Value = 1684.0731
$ws.Range("J138").Value = 3145.98
$ws.Range("K138").Value = 5052.219300000001
$ws.Range("L138").Value = 9437.940000000001
$ws.Range("M138").Value = 87.78069999999934
$ws.Range("N138").Value = -19717.94

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 16177.826
$ws.Range("I2").Value = 19354.736
$ws.Range("J2").Value = 1087.5
$ws.Range("K2").Value = 19354.736
$ws.Range("L2").Value = 1087.5
$ws.Range("M2").Value = -19241.736
$ws.Range("N2").Value = -1313.5

$ws.Range("H8").Value = 2502500
$ws.Range("I8").Value = 5000000
$ws.Range("J8").Value = 5000
$ws.Range("K8").Value = 5000000
$ws.Range("L8").Value = 5000
$ws.Range("M8").Value = -4999856
$ws.Range("N8").Value = -5288

$ws.Range("H10").Value = 45753.75
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 45753.75
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 45753.75
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -46093.75

$ws.Range("H12").Value = 4544.6665
$ws.Range("I12").Value = 8003
$ws.Range("J12").Value = 2815.5
$ws.Range("K12").Value = 8003
$ws.Range("L12").Value = 2815.5
$ws.Range("M12").Value = -7830
$ws.Range("N12").Value = -3161.5

$ws.Range("H13").Value = 2954.3333
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 2954.3333
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 2954.3333
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -3242.3333

$ws.Range("H43").Value = 14216.667
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 14216.667
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 14216.667
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -14842.667

$ws.Range("H116").Value = 16177.826
$ws.Range("I116").Value = 19354.736
$ws.Range("J116").Value = 1087.5
$ws.Range("K116").Value = 19354.736
$ws.Range("L116").Value = 1087.5
$ws.Range("M116").Value = -17060.736
$ws.Range("N116").Value = -5675.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 16177.826
$ws.Range("I3").Value = 19354.736
$ws.Range("J3").Value = 1087.5
$ws.Range("K3").Value = 19354.736
$ws.Range("L3").Value = 1087.5
$ws.Range("M3").Value = -19240.736
$ws.Range("N3").Value = -1315.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 5032.0557
$ws.Range("I122").Value = 7763.5
$ws.Range("K122").Value = 23290.5
$ws.Range("M122").Value = -20840.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2121054.2
$ws.Range("I5").Value = 2331426.8
$ws.Range("J5").Value = 1755670.1
$ws.Range("K5").Value = 6994280.399999999
$ws.Range("L5").Value = 5267010.300000001
$ws.Range("M5").Value = -6994168.399999999
$ws.Range("N5").Value = -5267234.300000001

$ws.Range("H7").Value = 66870064
$ws.Range("I7").Value = 77157740
$ws.Range("J7").Value = 140
$ws.Range("K7").Value = 231473220
$ws.Range("L7").Value = 420
$ws.Range("M7").Value = -231473108
$ws.Range("N7").Value = -644

$ws.Range("H8").Value = 145.8
$ws.Range("I8").Value = 145.8
$ws.Range("K8").Value = 437.4
$ws.Range("M8").Value = -298.4

$ws.Range("H80").Value = 6902.8
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 6902.8
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 20708.4
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -22580.4

$ws.Range("H83").Value = 6902.8
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 6902.8
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 62125.2
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -71485.20000000001

$ws.Range("H92").Value = 6098262
$ws.Range("J92").Value = 6098262
$ws.Range("L92").Value = 18294786
$ws.Range("N92").Value = -18297282

$ws.Range("H122").Value = 1274.091
$ws.Range("I122").Value = 298.125
$ws.Range("J122").Value = 1831.7858
$ws.Range("K122").Value = 2683.125
$ws.Range("L122").Value = 16486.0722
$ws.Range("M122").Value = -233.125
$ws.Range("N122").Value = -21386.0722

$ws.Range("H135").Value = 2121054.2
$ws.Range("I135").Value = 2331426.8
$ws.Range("J135").Value = 1755670.1
$ws.Range("K135").Value = 20982841.2
$ws.Range("L135").Value = 15801030.9
$ws.Range("M135").Value = -20980306.2
$ws.Range("N135").Value = -15806100.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 24181436
$ws.Range("I132").Value = 22511736
$ws.Range("J132").Value = 30303668
$ws.Range("K132").Value = 67535208
$ws.Range("L132").Value = 90911004
$ws.Range("M132").Value = -67532678
$ws.Range("N132").Value = -90916064

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
